$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores figures as literal text (e.g. "68.628.30" uses dots as
# thousands separators; values like "599.30" keep a significant trailing zero).
# Cells whose new text would otherwise be auto-parsed as a plain number by Excel
# (losing the trailing zero / being converted to a float) are switched to Text
# format first so the literal string is preserved exactly.
$numericLookingRows = @(5, 6, 8, 10, 12, 13, 14, 16, 20, 21, 22, 23, 25, 27, 29, 30, 32, 33, 34, 37, 38, 39, 40, 41, 42, 44, 47, 50, 51)
foreach ($r in $numericLookingRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "68.628.30"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "2.705.81"
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "599.30"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "160.82"
$ws.Range("E6").Value = "  +2.90%  "
$ws.Range("D8").Value = "0.543"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").Value = "2.705.15"
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  -4.68%  "
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "5.30"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "0.358"
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("D14").Value = "28.37"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").Value = "3.209.33"
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("D16").Value = "0.0000187"
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").Value = "68.574.41"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "2.768.18"
$ws.Range("E18").Value = "  +4.29%  "
$ws.Range("E19").Value = "  +4.04%  "
$ws.Range("D20").Value = "366.30"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").Value = "7.62"
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("D22").Value = "4.51"
$ws.Range("E22").Value = "  +2.79%  "
$ws.Range("D23").Value = "4.92"
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("E24").Value = "  +2.57%  "
$ws.Range("D25").Value = "75.15"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "10.15"
$ws.Range("E27").Value = "  +4.82%  "
$ws.Range("D28").Value = "2.818.78"
$ws.Range("D29").Value = "0.0000104"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("D30").Value = "581.64"
$ws.Range("E30").Value = "  +3.88%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").Value = "8.28"
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("D33").Value = "1.43"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").Value = "1.94"
$ws.Range("E34").Value = "  +4.74%  "
$ws.Range("E35").Value = "  +5.71%  "
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "20.15"
$ws.Range("E38").Value = "  +4.11%  "
$ws.Range("D39").Value = "161.99"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").Value = "0.381"
$ws.Range("E40").Value = "  +2.16%  "
$ws.Range("D41").Value = "1.89"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").Value = "5.42"
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").Value = "17.85"
$ws.Range("D45").Value = "0.0₆0319"
$ws.Range("E45").Value = "  -6.11%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "159.73"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("E48").Value = "  +4.89%  "
$ws.Range("E49").Value = "  +5.07%  "
$ws.Range("D50").Value = "0.606"
$ws.Range("E50").Value = "  +7.92%  "
$ws.Range("D51").Value = "22.21"
$ws.Range("E51").Value = "  +0.93%  "
